# Updated cryptos list with latest price/volume(1h) figures.
# Applies the per-row Price (column D) and Volume(1h) (column E) updates
# to Sheet1 of the workbook, matching the GitHub Actions data refresh.
#
# Column D values are stored as plain text in the source sheet (not
# numbers), so for any replacement text that looks numeric we force the
# cell's number format to Text ("@") first -- otherwise Excel would helpfully
# re-interpret e.g. "540.42" as the number 540.42 instead of keeping the
# literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.088.57'
$ws.Range("E2").Value = '  -4.59%  '
$ws.Range("D3").Value = '3.077.82'
$ws.Range("E3").Value = '  -4.84%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '540.42'
$ws.Range("E5").Value = '  -6.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.43'
$ws.Range("E6").Value = '  -11.95%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '3.072.40'
$ws.Range("E8").Value = '  -4.80%  '
$ws.Range("E9").Value = '  -4.83%  '
$ws.Range("E10").Value = '  -4.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.18'
$ws.Range("E11").Value = '  -12.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.459'
$ws.Range("E12").Value = '  -6.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000226'
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.39'
$ws.Range("E14").Value = '  -10.21%  '
$ws.Range("D15").Value = '3.533.59'
$ws.Range("E15").Value = '  -5.96%  '
$ws.Range("D16").Value = '63.016.17'
$ws.Range("E16").Value = '  -4.72%  '
$ws.Range("E17").Value = '  -3.48%  '
$ws.Range("D18").Value = '3.075.74'
$ws.Range("E18").Value = '  -4.85%  '
$ws.Range("E19").Value = '  -7.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '483.95'
$ws.Range("E20").Value = '  -11.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.28'
$ws.Range("E21").Value = '  -8.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.702'
$ws.Range("E22").Value = '  -5.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.16'
$ws.Range("E23").Value = '  -8.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.44'
$ws.Range("E24").Value = '  -3.55%  '
$ws.Range("E25").Value = '  -11.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("E27").Value = '  -8.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.16'
$ws.Range("E28").Value = '  -12.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.95'
$ws.Range("E30").Value = '  -6.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.89'
$ws.Range("E31").Value = '  -16.42%  '
$ws.Range("E32").Value = '  -6.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '58.93'
$ws.Range("E33").Value = '  +7.72%  '
$ws.Range("E34").Value = '  -12.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.95'
$ws.Range("E35").Value = '  -6.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.20'
$ws.Range("E36").Value = '  -7.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '465.50'
$ws.Range("E37").Value = '  -16.90%  '
$ws.Range("D38").Value = '3.125.65'
$ws.Range("E38").Value = '  -2.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0390'
$ws.Range("E39").Value = '  -14.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0787'
$ws.Range("E40").Value = '  -8.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.114'
$ws.Range("E41").Value = '  -12.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.04'
$ws.Range("E42").Value = '  -6.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.53'
$ws.Range("E43").Value = '  -12.71%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("E45").Value = '  -12.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.01'
$ws.Range("E46").Value = '  -13.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.32'
$ws.Range("E47").Value = '  -8.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '117.86'
$ws.Range("E48").Value = '  -5.26%  '
$ws.Range("E49").Value = '  -5.11%  '
$ws.Range("D50").Value = '0.0₃0509'
$ws.Range("E50").Value = '  -8.47%  '
$ws.Range("E51").Value = '  -9.98%  '
